$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.575.93'
$ws.Range('E2').Value = '  +4.01%  '
$ws.Range('D3').Value = '1.742.55'
$ws.Range('E3').Value = '  +4.22%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.000'
$ws.Range('E4').Value = '  +0.12%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '246.92'
$ws.Range('E5').Value = '  +3.10%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.001'
$ws.Range('E6').Value = '  +0.06%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4811'
$ws.Range('E7').Value = '  +0.88%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2690'
$ws.Range('E8').Value = '  +2.60%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06254'
$ws.Range('E9').Value = '  +1.13%  '
$ws.Range('D10').Value = '1.742.63'
$ws.Range('E10').Value = '  +4.24%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07121'
$ws.Range('E11').Value = '  +1.96%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '15.79'
$ws.Range('E12').Value = '  +6.24%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.6204'
$ws.Range('E13').Value = '  +5.10%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.502'
$ws.Range('E14').Value = '  +2.78%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '77.50'
$ws.Range('E15').Value = '  +2.77%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.001'
$ws.Range('D17').Value = '26.582.35'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '1.001'
$ws.Range('E18').Value = '  +0.11%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.000006893'
$ws.Range('E19').Value = '  +1.89%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '11.70'
$ws.Range('E20').Value = '  +2.16%  '
$ws.Range('D21').Value = '1.965.91'
$ws.Range('E21').Value = '  +4.14%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.635'
$ws.Range('E22').Value = '  +3.98%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '8.822'
$ws.Range('E23').Value = '  +0.09%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '5.340'
$ws.Range('E24').Value = '  +1.18%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '135.83'
$ws.Range('E25').Value = '  -0.78%  '
$ws.Range('E26').Value = '  +2.43%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.815'
$ws.Range('E27').Value = '  +4.86%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.437'
$ws.Range('E28').Value = '  +3.50%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '107.39'
$ws.Range('E29').Value = '  +2.60%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.011'
$ws.Range('E30').Value = '  +0.79%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.749'
$ws.Range('E31').Value = '  +3.07%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.04597'
$ws.Range('E33').Value = '  +7.16%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.619'
$ws.Range('E34').Value = '  -0.01%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.6438'
$ws.Range('E35').Value = '  +5.69%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.9968'
$ws.Range('E36').Value = '  +3.98%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.9463'
$ws.Range('E37').Value = '  +6.21%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '112.97'
$ws.Range('E38').Value = '  +17.37%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.995'
$ws.Range('E39').Value = '  +7.02%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.428'
$ws.Range('E40').Value = '  -6.40%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.003'
$ws.Range('E41').Value = '  +0.28%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.742'
$ws.Range('E42').Value = '  +16.48%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.01508'
$ws.Range('E43').Value = '  +1.69%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.3913'
$ws.Range('E44').Value = '  +3.96%  '
$ws.Range('E45').Value = '  +7.72%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '6.679'
$ws.Range('E46').Value = '  +7.17%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.05329'
$ws.Range('E47').Value = '  +1.15%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '7.933'
$ws.Range('E48').Value = '  +6.96%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '30.74'
$ws.Range('E49').Value = '  +2.66%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.273'
$ws.Range('E50').Value = '  +5.48%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.3450'
$ws.Range('E51').Value = '  +3.08%  '
